$wb = $excel.ActiveWorkbook

# The workbook's internal "width" (character units) is derived from the
# COM ColumnWidth property via width = (Round(ColumnWidth * 6) + 5) / 6,
# i.e. it only lands on multiples of 1/6. The target width from the repo
# (17.2159881591797) isn't reachable exactly through this engine, so we
# use the ColumnWidth input that rounds to the nearest achievable value
# (17.166666666666668).
$newColWidth = 16.333333333333332

# --- Overview sheet: Status (E/F) moved to "Ready for handoff", timestamp bumped ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-01 14:52:45"
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet: Status (C) moved to "Ready for handoff", Handoff Datetime (H) bumped ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-01 14:52:40"
$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet: Status (C) moved to "Ready for handoff", Handoff Datetime (H) bumped ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-01 14:52:45"
$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
